$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Gamma2F"

# Fix tiny floating point precision differences (last-ULP correction)
$ws.Range("C13").Value = 0.9902919605113798
$ws.Range("F13").Value = 0.9902919605113798
$ws.Range("C15").Value = 0.9888161938134855
$ws.Range("F15").Value = 0.9888161938134855
$ws.Range("L15").Value = 0.9927088830837937

# Add new row 16, copying formatting from row 15 (so the style/format of
# column A & B match exactly, without introducing new style entries)
$ws.Range("A15:M15").Copy($ws.Range("A16"))

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = $ws.Range("B15").Text

$ws.Range("C16").Value = 0.9997970327128718
$ws.Range("D16").Value = 0.9697660149721196
$ws.Range("E16").Value = 0.9997385364559017
$ws.Range("F16").Value = 0.9997970327128718
$ws.Range("G16").Value = 0.9699415614231334
$ws.Range("H16").Value = 1.001031999558077
$ws.Range("I16").Value = 0.9941468951873086
$ws.Range("J16").Value = 0.9697660149721196
$ws.Range("K16").Value = 0.9847522757140106
$ws.Range("L16").Value = 0.9922746542134413
$ws.Range("M16").Value = 0.9890703400515687
